# RBA v2.5 - Atualizacao da Tela
# Replace the placeholder "Tre"/"TRE" family tokens with "Qwer"/"QWER" family
# tokens, both in the document body and in the page header, preserving the
# exact case pattern of each individual occurrence (and the one transposed
# "Qewr" pair), in document order.

$d = $word.ActiveDocument

# --- 1. Main document body: single bold placeholder "TERE" -> "QWER" ---
$bodyRng = $d.Content
$bodyRng.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 0, $false, "QWER", 1) | Out-Null

# --- 2. Page header: ordered sequence of "Tre" family tokens ---
$targets      = @("TRE", "TERE", "Tre", "Tre", "Tre", "Tre", "Tre", "tre", "tre", "tre")
$replacements = @("QWER", "QWER", "Qwer", "Qwer", "Qewr", "Qewr", "Qwer", "qwer", "qwer", "qwer")

$hdr = $d.Sections.Item(1).Headers.Item(1)

$searchRng = $hdr.Range.Duplicate
$searchRng.Start = $hdr.Range.Start
$searchRng.End = $hdr.Range.End

for ($i = 0; $i -lt $targets.Count; $i++) {
    $found = $searchRng.Find.Execute($targets[$i], $true, $true, $false, $false, $false, $true, 0, $false, $replacements[$i], 1)
    if ($found) {
        $searchRng.Collapse(0)
        $searchRng.End = $hdr.Range.End
    }
}
